$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.991.40'
$ws.Range('E2').Value = '  +1.60%  '

# Row 3
$ws.Range('D3').Value = '2.454.67'
$ws.Range('E3').Value = '  +2.76%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.33'
$ws.Range('E5').Value = '  +1.28%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.92'
$ws.Range('E6').Value = '  +4.80%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.581'
$ws.Range('E8').Value = '  -0.54%  '

# Row 9
$ws.Range('D9').Value = '2.451.96'
$ws.Range('E9').Value = '  +2.57%  '

# Row 10
$ws.Range('E10').Value = '  +0.45%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.65'
$ws.Range('E11').Value = '  +0.16%  '

# Row 12
$ws.Range('E12').Value = '  +1.67%  '

# Row 13
$ws.Range('E13').Value = '  +0.19%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.71'
$ws.Range('E14').Value = '  +2.94%  '

# Row 15
$ws.Range('D15').Value = '2.892.39'
$ws.Range('E15').Value = '  +2.66%  '

# Row 16
$ws.Range('D16').Value = '62.867.61'
$ws.Range('E16').Value = '  +1.48%  '

# Row 17
$ws.Range('E17').Value = '  +1.07%  '

# Row 18
$ws.Range('D18').Value = '2.453.16'
$ws.Range('E18').Value = '  +1.36%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.25'
$ws.Range('E19').Value = '  +2.41%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '338.65'
$ws.Range('E20').Value = '  -2.81%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.24'
$ws.Range('E21').Value = '  +0.94%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.71'
$ws.Range('E22').Value = '  -1.48%  '

# Row 23
$ws.Range('E23').Value = '  -0.19%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.94'
$ws.Range('E24').Value = '  -0.07%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.171'
$ws.Range('E25').Value = '  -1.25%  '

# Row 26
$ws.Range('E26').Value = '  -0.03%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.49'
$ws.Range('E27').Value = '  +2.01%  '

# Row 28
$ws.Range('B28').Value = 'SuiNetwork'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.43'
$ws.Range('E28').Value = '  +5.24%  '

# Row 29
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.01'
$ws.Range('E29').Value = '  -3.90%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.74'
$ws.Range('E30').Value = '  +6.27%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.83'
$ws.Range('E31').Value = '  +2.22%  '

# Row 32
$ws.Range('D32').Value = '0.0₃0784'
$ws.Range('E32').Value = '  +2.84%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '176.79'
$ws.Range('E33').Value = '  +3.31%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.51'
$ws.Range('E34').Value = '  +8.41%  '

# Row 35
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '383.61'
$ws.Range('E35').Value = '  +10.48%  '

# Row 36
$ws.Range('B36').Value = 'PolygonEcosystemToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.397'
$ws.Range('E36').Value = '  +1.46%  '

# Row 37
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.72'
$ws.Range('E37').Value = '  +1.18%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.31'
$ws.Range('E39').Value = '  -3.73%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.11%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.71'
$ws.Range('E41').Value = '  +6.68%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '40.04'
$ws.Range('E42').Value = '  +2.69%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '148.51'
$ws.Range('E43').Value = '  +3.60%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.69'
$ws.Range('E44').Value = '  +0.03%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.39'
$ws.Range('E45').Value = '  +1.85%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.595'
$ws.Range('E46').Value = '  +2.00%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0958'
$ws.Range('E47').Value = '  -0.50%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0514'
$ws.Range('E48').Value = '  -0.47%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0228'
$ws.Range('E49').Value = '  +3.30%  '

# Row 50
$ws.Range('D50').Value = '0.0₆0232'
$ws.Range('E50').Value = '  +6.31%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.82'
$ws.Range('E51').Value = '  +0.87%  '
